# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that when the list of recorders ends with "System", "System" is moved to
# the front and the remaining order of the comma-separated list is reversed.
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#      "backup@backdoor.com, system, System" -> "System, system, backup@backdoor.com"
# Rows whose value does not end with "System" (e.g. single names, or lists
# ending in another name) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $parts = $text -split ", "
    $count = $parts.Length

    if ($count -gt 1 -and $parts[$count - 1] -eq "System") {
        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $newValue = [string]::Join(", ", $reversed)
        $cell.Value = $newValue
    }
}
